# Combine "Create own post(c)" / Read / Update / Delete own-post rows with
# the new buzzSpace "use case" column (C4:C7), renaming the B4 use-case
# label, and moving the active selection, per the commit:
#   "Creating of folder.  Combining Use Cases.  Also edited by Lerato."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "use case" descriptions added in column C for rows 4-7
$ws.Cells.Item(4, 3).Value = "User creates a post on the the buzzSpace. (l)"

# B4: "Create own post(c)" -> "create post"
$ws.Cells.Item(4, 2).Value = "create post"

$ws.Cells.Item(5, 3).Value = "user reads a post on the buzzSpace"
$ws.Cells.Item(6, 3).Value = "user edits and then update a post on the buzzSpace."
$ws.Cells.Item(7, 3).Value = "User can delete its own post."

# Move the saved view/selection (as captured when the edit was made)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C9").Select()
